$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/span"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$el = $wb.Worksheets.Item("Elements")

# Root "Extension" row no longer carries the ele-1/ext-1 constraint text
$el.Range("AI2").Value = ""

# Update ibm.com -> linuxforhealth.org URLs embedded in the extension slice type descriptions
$el.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/offset-begin}`n"
$el.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/offset-end}`n"
$el.Range("J7").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/covered-text}`n"
$el.Range("J8").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-confidence}`n"

# Fixed Value for Extension.url row
$el.Range("Q9").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/span"
